$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the start-of-credit dates in column A (rows 2-6)
$ws.Range("A2").Value = (Get-Date -Year 2022 -Month 10 -Day 2).Date
$ws.Range("A3").Value = (Get-Date -Year 2022 -Month 10 -Day 23).Date
$ws.Range("A4").Value = (Get-Date -Year 2022 -Month 10 -Day 24).Date
$ws.Range("A5").Value = (Get-Date -Year 2022 -Month 10 -Day 25).Date
$ws.Range("A6").Value = (Get-Date -Year 2022 -Month 10 -Day 26).Date

# Update the active selection/cell shown in the sheet view
$ws.Range("F5").Select()
